$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rolling-average sample sizes (column J = n) for several polls
$ws.Cells.Item(148, 10).Value = 1006
$ws.Cells.Item(150, 10).Value = 1108
$ws.Cells.Item(164, 10).Value = 1015
$ws.Cells.Item(175, 10).Value = 1057
$ws.Cells.Item(186, 10).Value = 1112

# Row 188: new cluster17 poll (2/1)
$ws.Cells.Item(188, 1).Value = 82
$ws.Cells.Item(188, 2).Value = 2022
$ws.Cells.Item(188, 3).Value = 22
$ws.Cells.Item(188, 4).Value = 1
$ws.Cells.Item(188, 5).Value = 27
$ws.Cells.Item(188, 6).Value = "cluster17"
$ws.Cells.Item(188, 7).Value = "online"
$ws.Cells.Item(188, 8).Value = "partially"
$ws.Cells.Item(188, 9).Value = 0
$ws.Cells.Item(188, 10).Value = 1997
$ws.Cells.Item(188, 11).Value = 1
$ws.Cells.Item(188, 12).Value = "T_0.5"
$ws.Cells.Item(188, 13).Value = 13
$ws.Cells.Item(188, 14).Value = 2
$ws.Cells.Item(188, 16).Value = 5
$ws.Cells.Item(188, 17).Value = 2
$ws.Cells.Item(188, 18).Value = 22.5
$ws.Cells.Item(188, 19).Value = 14
$ws.Cells.Item(188, 22).Value = 1
$ws.Cells.Item(188, 23).Value = 1.5
$ws.Cells.Item(188, 24).Value = 14.5
$ws.Cells.Item(188, 25).Value = 14.5
$ws.Cells.Item(188, 26).Value = 1
$ws.Cells.Item(188, 28).Value = 1
$ws.Cells.Item(188, 30).Value = 6

# Row 189: ifop rolling (1/31)
$ws.Cells.Item(189, 1).Value = 83
$ws.Cells.Item(189, 2).Value = 2022
$ws.Cells.Item(189, 3).Value = 22
$ws.Cells.Item(189, 4).Value = 1
$ws.Cells.Item(189, 5).Value = 31
$ws.Cells.Item(189, 6).Value = "ifop"
$ws.Cells.Item(189, 7).Value = "online"
$ws.Cells.Item(189, 8).Value = "included"
$ws.Cells.Item(189, 9).Value = 1
$ws.Cells.Item(189, 10).Value = 1000
$ws.Cells.Item(189, 11).Value = 0.5
$ws.Cells.Item(189, 12).Value = "T_0.5"
$ws.Cells.Item(189, 13).Value = 9.5
$ws.Cells.Item(189, 14).Value = 3.5
$ws.Cells.Item(189, 16).Value = 5
$ws.Cells.Item(189, 17).Value = 3.5
$ws.Cells.Item(189, 18).Value = 24
$ws.Cells.Item(189, 19).Value = 16
$ws.Cells.Item(189, 22).Value = 1
$ws.Cells.Item(189, 23).Value = 1.5
$ws.Cells.Item(189, 24).Value = 17.5
$ws.Cells.Item(189, 25).Value = 14
$ws.Cells.Item(189, 30).Value = 4
$ws.Cells.Item(189, 31).Value = "T_0.5"

# Row 190: opinionway rolling (2/1)
$ws.Cells.Item(190, 1).Value = 84
$ws.Cells.Item(190, 2).Value = 2022
$ws.Cells.Item(190, 3).Value = 22
$ws.Cells.Item(190, 4).Value = 1
$ws.Cells.Item(190, 5).Value = 29
$ws.Cells.Item(190, 6).Value = "opinionway"
$ws.Cells.Item(190, 7).Value = "online"
$ws.Cells.Item(190, 8).Value = "partially"
$ws.Cells.Item(190, 9).Value = 1
$ws.Cells.Item(190, 10).Value = 800
$ws.Cells.Item(190, 11).Value = "T_1"
$ws.Cells.Item(190, 12).Value = "T_1"
$ws.Cells.Item(190, 13).Value = 10
$ws.Cells.Item(190, 14).Value = 3
$ws.Cells.Item(190, 16).Value = 5
$ws.Cells.Item(190, 17).Value = 3
$ws.Cells.Item(190, 18).Value = 24
$ws.Cells.Item(190, 19).Value = 16
$ws.Cells.Item(190, 22).Value = 1
$ws.Cells.Item(190, 23).Value = 2
$ws.Cells.Item(190, 24).Value = 17
$ws.Cells.Item(190, 25).Value = 14
$ws.Cells.Item(190, 30).Value = 5

# Update sheet view: frozen pane position and active selection
$ws.Activate()
$appWin = $excel.ActiveWindow
$appWin.ScrollColumn = 12
$appWin.SplitRow = 1
$appWin.FreezePanes = $true
$ws.Range("A179").Select()
$excel.ActiveWindow.ScrollRow = 179
$ws.Range("U189").Select()
